$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the style of the existing header (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting from H1 (bold header style) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in column I (I0) and column J (IF) values for rows 2-26
$iValues = @{
    2 = 3
    3 = 1
    4 = 1
    5 = 1
    6 = 1
    7 = 1
    8 = 1
    9 = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 1
    26 = 1
}

$jValues = @{
    2 = 5
    3 = 4
    4 = 5
    5 = 6
    6 = 5
    7 = 6
    8 = 5
    9 = 5
    10 = 6
    11 = 5
    12 = 7
    13 = 6
    14 = 6
    15 = 5
    16 = 7
    17 = 5
    18 = 5
    19 = 5
    20 = 3
    21 = 6
    22 = 9
    23 = 6
    24 = 4
    25 = 3
    26 = 2
}

for ($row = 2; $row -le 26; $row++) {
    $ws.Cells.Item($row, 9).Value = $iValues[$row]
    $ws.Cells.Item($row, 10).Value = $jValues[$row]
}
